$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final (post-edit) values for columns C:H, rows 2..21.
# Row 2 holds new sensor readings; rows 3..21 are the old rows 2..20
# shifted down by one (the old last row, 22, is dropped).
$data = @(
    @(-3.555192089080811,4.907798504829406,-2.923101136088372,-0.0018325957935303,-0.0296269636601209,-0.0087048299610614),
    @(-3.523949909210205,4.786228704452514,-1.791911107301711,0.0019853119738399,0.0355829000473022,0.0200058370828628),
    @(-4.734984159469604,4.268091917037964,-0.5654808729887018,0.0106901414692401,0.08765916526317589,0.0525344125926494),
    @(-5.463105487823487,3.789239883422852,0.125382423400879,-0.0532979927957057,0.3320052623748779,0.1557706445455551),
    @(-6.153488552570344,2.863918662071227,-0.2511623546481139,-0.1586722433567047,0.3882048726081848,-0.0940732508897781),
    @(-7.365730333328246,1.174339866638185,0.3271868914365764,-0.2698497176170349,0.6501133441925049,-0.4689917862415313),
    @(-7.554701018333435,1.320136770606041,0.9918349981307982,0.0048869219608604,0.8135197758674622,-1.29075825214386),
    @(-6.541906356811527,3.541345179080956,-0.3953665494918761,0.8275696635246277,-2.131614208221436,-2.630385637283325),
    @(41.03047697544085,6.894529008865351,-6.942315888404829,1.503492116928101,-4.445571899414063,-2.242486238479614),
    @(64.96824674606313,6.919658172130577,-9.210757869482023,-1.950034618377685,3.770718336105346,0.2194533348083496),
    @(-6.831269806623445,3.269432669878014,-0.1760272979736253,6.251442432403564,-2.222022294998169,-3.614794969558716),
    @(-2.175523471832268,5.449279594421388,1.586600971221922,0.8747590184211731,2.045024156570435,1.580461144447327),
    @(2.27016156911849,5.673981070518497,0.4227319359779442,-0.0520762614905834,-0.5412266254425049,0.422871470451355),
    @(4.514346599578862,5.118150138854979,-1.915704894065858,-0.2570215463638305,-7.184691429138184,2.630538463592529),
    @(-6.036725759506252,9.008556652069101,-0.9792453408241246,2.87503719329834,-1.23394775390625,-0.1574505120515823),
    @(-1.178059291839593,4.555799674987787,3.518247509002693,-0.157145082950592,1.678810358047485,1.309847831726074),
    @(-1.642422831058502,6.837077736854553,3.844039088487625,0.982424020767212,1.2322678565979,-0.3689626157283783),
    @(0.4585734605789122,6.029543757438664,3.976011931896209,0.0326812900602817,0.007177666760981,-0.1241583600640297),
    @(0.8665444850921629,5.183062970638275,4.418135178089142,-0.00167987938039,0.0308486949652433,0.0305432621389627),
    @(0.2307996749877924,4.777379417419433,4.424502086639404,0.012980886735022,-0.06704246252775189,-0.0529925599694252)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    for ($j = 0; $j -lt $vals.Length; $j++) {
        $col = $j + 3
        $ws.Cells.Item($row, $col).Value = $vals[$j]
    }
}

# Remove the now-obsolete last row (old row 22), shrinking the used range to A1:H21.
$ws.Rows.Item(22).Delete()
